# Regenerate orders with updated distance/size codes.
# The experiment's distance and size condition labels changed:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# These codes appear embedded inside many strings throughout the sheet
# (condition names, filenames, and the standalone Distance/Size columns),
# so we perform a whole-sheet text substitution for each token.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.Cells

$cells.Replace("D51", "D55")
$cells.Replace("D64", "D69")
$cells.Replace("D80", "D86")
$cells.Replace("S30", "S31")
